$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00006486019690155054
$ws.Range("C2").Value = 0.0001537489499301437
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.6530110982390698
